# unitTest_base_macro3.xlsx - add new base command `assertMatch(text,regex)`
# and new external command `openFile(filePath)` to the hidden '#system' sheet,
# and remove the obsolete `tn.5250` command column (and its defined name).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) "target" column (A): the category list loses the "tn.5250" row (A27).
#    Shift A28:A33 up into A27:A32, then blank out A33.
# ---------------------------------------------------------------------------
for ($r = 27; $r -le 32; $r++) {
    $srcVal = $ws.Cells.Item($r + 1, 1).Value()
    $ws.Cells.Item($r, 1).Value = $srcVal
}
$ws.Cells.Item(33, 1).Value = ""

# ---------------------------------------------------------------------------
# 2) "base" column (F): insert new command assertMatch(text,regex) at F11,
#    keeping alphabetical order (after assertEqual, before assertNotContain).
#    Shift F11:F44 down into F12:F45, then set the new F11 value.
# ---------------------------------------------------------------------------
for ($r = 45; $r -ge 12; $r--) {
    $srcVal = $ws.Cells.Item($r - 1, 6).Value()
    $ws.Cells.Item($r, 6).Value = $srcVal
}
$ws.Cells.Item(11, 6).Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------------
# 3) "external" column (J): insert new command openFile(filePath) at J2,
#    keeping alphabetical order (before runJUnit).
#    Shift J2:J6 down into J3:J7, then set the new J2 value.
# ---------------------------------------------------------------------------
for ($r = 7; $r -ge 3; $r--) {
    $srcVal = $ws.Cells.Item($r - 1, 10).Value()
    $ws.Cells.Item($r, 10).Value = $srcVal
}
$ws.Cells.Item(2, 10).Value = "openFile(filePath)"

# ---------------------------------------------------------------------------
# 4) "tn.5250" column (AA) is removed entirely; columns AB:AG shift left
#    one position to AA:AF, for every row (1-151).
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 151; $r++) {
    for ($c = 27; $c -le 32; $c++) {
        $srcVal = $ws.Cells.Item($r, $c + 1).Value()
        $ws.Cells.Item($r, $c).Value = $srcVal
    }
    $ws.Cells.Item($r, 33).Value = ""
}

# ---------------------------------------------------------------------------
# 5) Update the workbook-level defined names to reflect the new ranges.
# ---------------------------------------------------------------------------
$names = $wb.Names

$names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$45"
$names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$7"
$names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$151"
$names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"

$names.Item("tn.5250").Delete()
